# Add new columns I ("I0") and J ("IF") to Sheet1, mirroring the existing
# header style used by the other header cells (e.g. H1), and fill in the
# per-row numeric values for rows 2-31.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Match the style already used by the other header cells in row 1 by
# copying the formatting from an existing header cell (H1) before
# writing the new header text.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-31: column I ("I0") then column J ("IF")
$values = @{
    2  = @(8, 9)
    3  = @(8, 8)
    4  = @(7, 8)
    5  = @(7, 7)
    6  = @(5, 6)
    7  = @(5, 6)
    8  = @(7, 7)
    9  = @(5, 6)
    10 = @(9, 9)
    11 = @(8, 8)
    12 = @(5, 6)
    13 = @(6, 6)
    14 = @(8, 8)
    15 = @(8, 8)
    16 = @(5, 6)
    17 = @(9, 9)
    18 = @(4, 5)
    19 = @(6, 6)
    20 = @(8, 8)
    21 = @(6, 6)
    22 = @(7, 7)
    23 = @(5, 5)
    24 = @(7, 7)
    25 = @(7, 7)
    26 = @(8, 8)
    27 = @(3, 3)
    28 = @(5, 5)
    29 = @(8, 8)
    30 = @(8, 8)
    31 = @(7, 7)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
